$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "37.813.36"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.35%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.105.38"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.30%  "

$ws.Range("E4").Value = "  +0.07%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "234.59"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.02%  "

$ws.Range("E6").Value = "  +0.97%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "57.96"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.08%  "

$ws.Range("E8").Value = "  +0.07%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.391"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.02%  "

$ws.Range("E10").Value = "  +2.19%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "2.414.53"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.44%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "14.46"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.24%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "21.34"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.94%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.782"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.25%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.22"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.41%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.087.83"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.21%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "37.759.60"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.33%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.22"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.85%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "70.25"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.50%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0822"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.24%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "227.34"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.96%  "

$ws.Range("E23").Value = "  -0.03%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.42"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.52%  "

$ws.Range("E25").Value = "  +0.88%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "168.81"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.78%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.95"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.54%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.133"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +4.08%  "

$ws.Range("E29").Value = "  -3.25%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "19.42"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.85%  "

$ws.Range("E31").Value = "  +1.17%  "

$ws.Range("E32").Value = "  +3.48%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0624"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.41%  "

$ws.Range("E34").Value = "  +0.64%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.58"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.83%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.46"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +5.56%  "

$ws.Range("E37").Value = "  +3.52%  "

$ws.Range("E38").Value = "  +0.02%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.44"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -7.02%  "

$ws.Range("E40").Value = "  +6.53%  "

$ws.Range("E41").Value = "  -0.18%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "96.69"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.12%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.467.77"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.62%  "

$ws.Range("E44").Value = "  +1.39%  "

$ws.Range("E45").Value = "  +0.89%  "

$ws.Range("E46").Value = "  -11.69%  "

$ws.Range("E47").Value = "  +2.95%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.04"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.05%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "15.37"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.94%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.29"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.01%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.300.66"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.40%  "
